# Updated symbol list on Tue Feb 14 06:58:54 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures
# for the cryptocurrency exchange-token rows on the active sheet. Both
# columns hold plain text (e.g. "292.06", "-7.00%") rather than numbers,
# so each target cell's NumberFormat is forced to "@" (Text) before the
# new value is written — this keeps Excel from auto-coercing the digits
# or the trailing "%" into a numeric/percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "292.06"
    "E2" = "-7.00%"
    "D3" = "40.48"
    "E3" = "-0.87%"
    "D4" = "5.033"
    "E4" = "-2.48%"
    "D5" = "0.07324"
    "E5" = "-3.54%"
    "E6" = "-9.05%"
    "D7" = "0.9305"
    "E7" = "0.08%"
    "D8" = "2.332"
    "E8" = "-3.80%"
    "D9" = "0.1177"
    "E9" = "-1.95%"
    "D10" = "0.1739"
    "E10" = "-4.37%"
    "D11" = "0.04336"
    "E11" = "4.65%"
    "D12" = "0.08686"
    "E12" = "-3.99%"
    "E13" = "0.04%"
    "D14" = "0.001266"
    "E14" = "-1.70%"
    "D15" = "0.005854"
    "E15" = "0.31%"
    "E16" = "0.04%"
    "D17" = "4.279"
    "E17" = "-1.10%"
    "E18" = "-1.71%"
    "D19" = "7.976"
    "E19" = "4.63%"
    "D20" = "0.1400"
    "E20" = "4.32%"
    "D21" = "0.2742"
    "E21" = "-3.45%"
    "D22" = "0.03934"
    "E22" = "-1.76%"
    "E23" = "-1.49%"
    "D24" = "0.003788"
    "E24" = "-4.64%"
    "E25" = "-5.25%"
    "D26" = "0.0003724"
    "D38" = "0.02282"
    "E38" = "-5.36%"
    "D39" = "0.05038"
    "E39" = "-2.51%"
    "D40" = "0.006287"
    "E40" = "90.27%"
    "D41" = "0.007686"
    "E41" = "-0.64%"
    "E42" = "-0.80%"
    "E43" = "-3.58%"
    "D44" = "0.008274"
    "E44" = "-3.49%"
    "D45" = "0.2915"
    "E45" = "-13.95%"
    "D46" = "0.00006278"
    "E46" = "-4.88%"
    "E47" = "-0.02%"
    "D48" = "0.03196"
    "E48" = "-88.39%"
    "E49" = "-0.02%"
    "E50" = "-0.02%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
